$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1282158811859517
$ws.Range("D2").Value = 0.08878388640228341
$ws.Range("E2").Value = 0.1065813741126328
$ws.Range("F2").Value = 2.118513577096095
$ws.Range("G2").Value = 0.002497283076363639
$ws.Range("I2").Value = 1.509516333546614
$ws.Range("J2").Value = 0.1455084876916897
$ws.Range("L2").Value = 0.1547937879303847
$ws.Range("M2").Value = 1.489323724784199
$ws.Range("N2").Value = 1.814412066740545
$ws.Range("O2").Value = 5.884379728857311

$ws.Range("C3").Value = 0.1291791960709752
$ws.Range("D3").Value = 0.08929637871369467
$ws.Range("E3").Value = 0.1077315201928244
$ws.Range("F3").Value = 2.106293471411519
$ws.Range("G3").Value = 0.002501595174405513
$ws.Range("I3").Value = 1.499915050636943
$ws.Range("J3").Value = 0.1474638540712763
$ws.Range("L3").Value = 0.1564256208648072
$ws.Range("M3").Value = 1.37519768558542
$ws.Range("N3").Value = 1.675437757530716
$ws.Range("O3").Value = 5.834295394609342

$ws.Range("C4").Value = 0.1298165073425359
$ws.Range("D4").Value = 0.08963032639650526
$ws.Range("E4").Value = 0.1084756106854649
$ws.Range("F4").Value = 2.10001937566426
$ws.Range("G4").Value = 0.002504385820531424
$ws.Range("I4").Value = 1.494905133959001
$ws.Range("J4").Value = 0.1487277537139693
$ws.Range("L4").Value = 0.1574810136276206
$ws.Range("M4").Value = 1.305026421810354
$ws.Range("N4").Value = 1.590341961381142
$ws.Range("O4").Value = 5.807067093090922

$ws.Range("C5").Value = 0.1300877464390098
$ws.Range("D5").Value = 0.08977126629515819
$ws.Range("E5").Value = 0.1087883638078747
$ws.Range("F5").Value = 2.097771241142681
$ws.Range("G5").Value = 0.002505559105454628
$ws.Range("I5").Value = 1.493085821731498
$ws.Range("J5").Value = 0.1492587082314114
$ws.Range("L5").Value = 0.1579245387751005
$ws.Range("M5").Value = 1.276409384577647
$ws.Range("N5").Value = 1.555727969744765
$ws.Range("O5").Value = 5.796855853396437

$ws.Range("C6").Value = 0.1301334819832221
$ws.Range("D6").Value = 0.08979496259230402
$ws.Range("E6").Value = 0.1088408718540131
$ws.Range("F6").Value = 2.097416566219167
$ws.Range("G6").Value = 0.002505756110593084
$ws.Range("I6").Value = 1.492797140750547
$ws.Range("J6").Value = 0.1493478331267053
$ws.Range("L6").Value = 0.157998997949246
$ws.Range("M6").Value = 1.271656321247619
$ws.Range("N6").Value = 1.549984284947584
$ws.Range("O6").Value = 5.795213662573644

$ws.Range("C7").Value = 0.1298201186782357
$ws.Range("D7").Value = 0.08963220750051426
$ws.Range("E7").Value = 0.1084797899950176
$ws.Range("F7").Value = 2.099987807596804
$ws.Range("G7").Value = 0.002504401497534206
$ws.Range("I7").Value = 1.49487969858523
$ws.Range("J7").Value = 0.1487348499684022
$ws.Range("L7").Value = 0.1574869407213599
$ws.Range("M7").Value = 1.304640566469615
$ws.Range("N7").Value = 1.589874883414069
$ws.Range("O7").Value = 5.806925801447363

$ws.Range("C8").Value = 0.1285385225888334
$ws.Range("D8").Value = 0.08895659713904536
$ws.Range("E8").Value = 0.1069700798100086
$ws.Range("F8").Value = 2.114044699569064
$ws.Range("G8").Value = 0.00249874027736834
$ws.Range("I8").Value = 1.506021847022666
$ws.Range("J8").Value = 0.1461695541679715
$ws.Range("L8").Value = 0.1553453560144487
$ws.Range("M8").Value = 1.449994935321982
$ws.Range("N8").Value = 1.766447378977404
$ws.Range("O8").Value = 5.86637834463744

$ws.Range("C9").Value = 0.1263887608402463
$ws.Range("D9").Value = 0.08778436322862193
$ws.Range("E9").Value = 0.1043101448843785
$ws.Range("F9").Value = 2.151387672161107
$ws.Range("G9").Value = 0.002488768041149199
$ws.Range("I9").Value = 1.534915181455617
$ws.Range("J9").Value = 0.141641673551899
$ws.Range("L9").Value = 0.1515694178027012
$ws.Range("M9").Value = 1.73415175492147
$ws.Range("N9").Value = 2.114407389216467
$ws.Range("O9").Value = 6.011007824062858

$ws.Range("C10").Value = 0.1250306052860068
$ws.Range("D10").Value = 0.08701573927993689
$ws.Range("E10").Value = 0.1025389602742741
$ws.Range("F10").Value = 2.184825430925684
$ws.Range("G10").Value = 0.00248212261082123
$ws.Range("I10").Value = 1.560468331020388
$ws.Range("J10").Value = 0.1386219370140054
$ws.Range("L10").Value = 0.1490530025230452
$ws.Range("M10").Value = 1.942256951500383
$ws.Range("N10").Value = 2.370901038614818
$ws.Range("O10").Value = 6.134501091156494

$ws.Range("C11").Value = 0.124460733829693
$ws.Range("D11").Value = 0.08668609057558818
$ws.Range("E11").Value = 0.1017729132946903
$ws.Range("F11").Value = 2.201349785146888
$ws.Range("G11").Value = 0.00247924577941853
$ws.Range("I11").Value = 1.573039404167503
$ws.Range("J11").Value = 0.1373149188129448
$ws.Range("L11").Value = 0.1479640768273143
$ws.Range("M11").Value = 2.036759658221371
$ws.Range("N11").Value = 2.487732864659222
$ws.Range("O11").Value = 6.194454942899995

$ws.Range("C12").Value = 0.1242518312997234
$ws.Range("D12").Value = 0.08656413082183256
$ws.Range("E12").Value = 0.1014885353823209
$ws.Range("F12").Value = 2.207796635472619
$ws.Range("G12").Value = 0.002478177303627155
$ws.Range("I12").Value = 1.577936385353112
$ws.Range("J12").Value = 0.1368295853630683
$ws.Range("L12").Value = 0.1475597477616652
$ws.Range("M12").Value = 2.07251904454148
$ws.Range("N12").Value = 2.531992182504212
$ws.Range("O12").Value = 6.217703074897713

$ws.Range("C13").Value = 0.1242965155404434
$ws.Range("D13").Value = 0.08659026942862491
$ws.Range("E13").Value = 0.1015495273742395
$ws.Range("F13").Value = 2.206399756413219
$ws.Range("G13").Value = 0.002478406490543528
$ws.Range("I13").Value = 1.576875651001501
$ws.Range("J13").Value = 0.1369336830848487
$ws.Range("L13").Value = 0.1476464704309883
$ws.Range("M13").Value = 2.064818856875831
$ws.Range("N13").Value = 2.522459427520289
$ws.Range("O13").Value = 6.212671908155016

$ws.Range("C14").Value = 0.1244434091049733
$ws.Range("D14").Value = 0.08667599936223347
$ws.Range("E14").Value = 0.1017494029526871
$ws.Range("F14").Value = 2.201876371490528
$ws.Range("G14").Value = 0.002479157456697927
$ws.Range("I14").Value = 1.573439541795125
$ws.Range("J14").Value = 0.1372747975150697
$ws.Range("L14").Value = 0.147930651635475
$ws.Range("M14").Value = 2.039702157171973
$ws.Range("N14").Value = 2.491373775439115
$ws.Range("O14").Value = 6.196356648932408

$ws.Range("C15").Value = 0.1245342836089947
$ws.Range("D15").Value = 0.08672888516428756
$ws.Range("E15").Value = 0.101872575867525
$ws.Range("F15").Value = 2.199130355036942
$ws.Range("G15").Value = 0.002479620165885882
$ws.Range("I15").Value = 1.571352626624503
$ws.Range("J15").Value = 0.1374849914780718
$ws.Range("L15").Value = 0.1481057656902429
$ws.Range("M15").Value = 2.024313883985684
$ws.Range("N15").Value = 2.472335090904153
$ws.Range("O15").Value = 6.186434101739337

$ws.Range("C16").Value = 0.1250688127057487
$ws.Range("D16").Value = 0.08703768464527739
$ws.Range("E16").Value = 0.1025898217884214
$ws.Range("F16").Value = 2.183772002723828
$ws.Range("G16").Value = 0.002482313551772035
$ws.Range("I16").Value = 1.559665871492726
$ws.Range("J16").Value = 0.1387086969496163
$ws.Range("L16").Value = 0.1491252892031003
$ws.Range("M16").Value = 1.936077396231639
$ws.Range("N16").Value = 2.363268497672323
$ws.Range("O16").Value = 6.130659113964271

$ws.Range("C17").Value = 0.1254090118869549
$ws.Range("D17").Value = 0.08723224242506333
$ws.Range("E17").Value = 0.1030399920299537
$ws.Range("F17").Value = 2.174686940557933
$ws.Range("G17").Value = 0.002484003227720133
$ws.Range("I17").Value = 1.552739238071879
$ws.Range("J17").Value = 0.1394764905791259
$ws.Range("L17").Value = 0.1497650258105967
$ws.Range("M17").Value = 1.881902699558992
$ws.Range("N17").Value = 2.296395524080367
$ws.Range("O17").Value = 6.097411475327362

$ws.Range("C18").Value = 0.1256091997782995
$ws.Range("D18").Value = 0.08734602992374896
$ws.Range("E18").Value = 0.1033026521831952
$ws.Range("F18").Value = 2.169585036895199
$ws.Range("G18").Value = 0.00248498885255467
$ws.Range("I18").Value = 1.548844314225406
$ws.Range("J18").Value = 0.1399243798450573
$ws.Range("L18").Value = 0.1501382380038088
$ws.Range("M18").Value = 1.850727454841788
$ws.Range("N18").Value = 2.257946470643219
$ws.Range("O18").Value = 6.078643626572671

$ws.Range("C19").Value = 0.1256777553358681
$ws.Range("D19").Value = 0.08738487998541267
$ws.Range("E19").Value = 0.1033922256667504
$ws.Range("F19").Value = 2.167878828100498
$ws.Range("G19").Value = 0.002485324936607205
$ws.Range("I19").Value = 1.547540848396579
$ws.Range("J19").Value = 0.1400771045855687
$ws.Range("L19").Value = 0.1502655034178204
$ws.Range("M19").Value = 1.840169499412497
$ws.Range("N19").Value = 2.244930901807322
$ws.Range("O19").Value = 6.072350123185345

$ws.Range("C20").Value = 0.1253723298976439
$ws.Range("D20").Value = 0.08721133657027824
$ws.Range("E20").Value = 0.1029916841538812
$ws.Range("F20").Value = 2.175641266348222
$ws.Range("G20").Value = 0.002483821934739306
$ws.Range("I20").Value = 1.553467366159353
$ws.Range("J20").Value = 0.1393941081021914
$ws.Range("L20").Value = 0.1496963811407301
$ws.Range("M20").Value = 1.887671302993695
$ws.Range("N20").Value = 2.30351279459552
$ws.Range("O20").Value = 6.100913955973624

$ws.Range("C21").Value = 0.1244000757734689
$ws.Range("D21").Value = 0.08665074052422028
$ws.Range("E21").Value = 0.101690539732038
$ws.Range("F21").Value = 2.203199853413224
$ws.Range("G21").Value = 0.002478936312774372
$ws.Range("I21").Value = 1.574445099875007
$ws.Range("J21").Value = 0.1371743430889047
$ws.Range("L21").Value = 0.1478469629802444
$ws.Range("M21").Value = 2.04708028981338
$ws.Range("N21").Value = 2.500503939066562
$ws.Range("O21").Value = 6.201134032083473

$ws.Range("C22").Value = 0.1238048428051641
$ws.Range("D22").Value = 0.08630109092154825
$ws.Range("E22").Value = 0.1008734352113634
$ws.Range("F22").Value = 2.222315434104658
$ws.Range("G22").Value = 0.002475865153594371
$ws.Range("I22").Value = 1.588951610996958
$ws.Range("J22").Value = 0.1357795978976215
$ws.Range("L22").Value = 0.1466850277778882
$ws.Range("M22").Value = 2.151106073123856
$ws.Range("N22").Value = 2.629349973736907
$ws.Range("O22").Value = 6.269810693701857

$ws.Range("C23").Value = 0.1241188527117423
$ws.Range("D23").Value = 0.08648617620309684
$ws.Range("E23").Value = 0.1013064943803852
$ws.Range("F23").Value = 2.212011838959157
$ws.Range("G23").Value = 0.002477493172220034
$ws.Range("I23").Value = 1.581136201177387
$ws.Range("J23").Value = 0.1365188696318986
$ws.Range("L23").Value = 0.1473008956557766
$ws.Range("M23").Value = 2.095600910027997
$ws.Range("N23").Value = 2.56057454390276
$ws.Range("O23").Value = 6.232865332144684

$ws.Range("C24").Value = 0.1253888994876
$ws.Range("D24").Value = 0.08722078208362927
$ws.Range("E24").Value = 0.1030135121486409
$ws.Range("F24").Value = 2.17520943811391
$ws.Range("G24").Value = 0.002483903852892531
$ws.Range("I24").Value = 1.553137907583718
$ws.Range("J24").Value = 0.1394313330538697
$ws.Range("L24").Value = 0.149727398515795
$ws.Range("M24").Value = 1.885063409008268
$ws.Range("N24").Value = 2.300295085279402
$ws.Range("O24").Value = 6.099329404796322

$ws.Range("C25").Value = 0.1269314495910727
$ws.Range("D25").Value = 0.08808518794043918
$ws.Range("E25").Value = 0.1049975657000375
$ws.Range("F25").Value = 2.140234422425578
$ws.Range("G25").Value = 0.002491345646585021
$ws.Range("I25").Value = 1.526341490042796
$ws.Range("J25").Value = 0.1428127252156757
$ws.Range("L25").Value = 0.1525456046951472
$ws.Range("M25").Value = 1.657387476139931
$ws.Range("N25").Value = 2.020112649522616
$ws.Range("O25").Value = 5.968865299716242

Write-Host "Applied 380 kV case updates"